# "Generate Report for Handoff"
#
# The localization status report is being (re)generated for a fresh
# handoff: the per-language status that used to read "Handed back: in
# sync with en-US" is now "Ready for handoff", and the associated
# generation/handoff timestamps are refreshed.
#
#   Overview!E2, F2  (zh-cn / de-de status)   -> "Ready for handoff"
#   Overview!G2      (Latest HO Xliff Gen.)   -> 2016-09-04 03:04:11
#   zh-cn!C2         (Status)                 -> "Ready for handoff"
#   zh-cn!H2         (Latest Handoff Datetime) -> 2016-09-04 03:04:05
#   de-de!C2         (Status)                 -> "Ready for handoff"
#   de-de!H2         (Latest Handoff Datetime) -> 2016-09-04 03:04:11
#
# The "Status"/language columns also get narrower (their header text is
# shorter now than the old "Handed back: in sync with en-US" phrase
# needed).

$wb = $excel.ActiveWorkbook

$statusText = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("G2").Value = "2016-09-04 03:04:11"
$wsOverview.Range("E1:F1").ColumnWidth = 16.3333333333333

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("H2").Value = "2016-09-04 03:04:05"
$wsZhCn.Range("C1").ColumnWidth = 16.3333333333333

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("H2").Value = "2016-09-04 03:04:11"
$wsDeDe.Range("C1").ColumnWidth = 16.3333333333333
